$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row at position 59, pushing the existing rows 59-64 down to 60-65.
$ws.Rows.Item(59).Insert()

# Populate the newly inserted row 59 with the new weekly data point.
$ws.Cells.Item(59, 1).Value = 9
$ws.Cells.Item(59, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(59, 3).Value = "Metropolitana"
$ws.Cells.Item(59, 4).Value = 44461
$ws.Cells.Item(59, 4).NumberFormat = $ws.Cells.Item(60, 4).NumberFormat
$ws.Cells.Item(59, 5).Value = 13
$ws.Cells.Item(59, 6).Value = 100112005
$ws.Cells.Item(59, 7).Value = "Puerro"
$ws.Cells.Item(59, 8).Value = "Sin especificar"
$ws.Cells.Item(59, 9).Value = "Primera"
$ws.Cells.Item(59, 10).Value = 79
$ws.Cells.Item(59, 11).Value = 7000
$ws.Cells.Item(59, 12).Value = 8000
$ws.Cells.Item(59, 13).Value = 7494
$ws.Cells.Item(59, 14).Value = "$/paquete 20 unidades"
$ws.Cells.Item(59, 15).Value = "Provincia de Chacabuco"
$ws.Cells.Item(59, 16).Value = 375
$ws.Cells.Item(59, 17).Value = 20
$ws.Cells.Item(59, 18).Value = "Hortaliza"
